$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 305.43452465743718
$ws.Range("C2").Value = 558.17631017886845
$ws.Range("D2").Value = 319.5879819441991
$ws.Range("E2").Value = 485.99970020148959

$ws.Range("B3").Value = 332.52661422423711
$ws.Range("C3").Value = 353.48246798397503
$ws.Range("D3").Value = 294.56426872549542
$ws.Range("E3").Value = 302.71437874738388

$ws.Range("B1:E3").Select()
